$d = $word.ActiveDocument

# Locate the paragraph that holds the "vie DK" hit-box variables
# (var vieDKX / vieDKY / vieDKL / vieDKH). This paragraph currently ends
# right after "var vieDKH = [33,33,33,33] ;" with the _GoBack bookmark.
$range = $d.Content
$found = $range.Find.Execute("var vieDKX = [172,178,184,189] ;", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'vieDKX' paragraph"
}

$para = $range.Paragraphs(1)
$prange = $para.Range

# Replace that single paragraph with three paragraphs:
#  1) the original vieDK* lines (now with a bold paragraph mark), three
#     extra line breaks, and a new bold "Score+vie :" heading;
#  2) the new vie (heart/life icon) sprite coordinates block;
#  3) the new Score sprite coordinates block, keeping the pre-existing
#     _GoBack bookmark positioned between "ScoreX= [188,203,219,235" and
#     the closing "]".
# An extra blank paragraph is appended to preserve the blank paragraph
# that originally followed this block.
$bodyXml = '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:t>var vieDKX = [172,178,184,189] ;</w:t></w:r><w:r><w:br/><w:t>var vieDKY = [51,51,51,51] ;</w:t></w:r><w:r><w:br/><w:t>var vieDKL = [10,6,6,7] ;</w:t></w:r><w:r><w:br/><w:t>var vieDKH = [33,33,33,33] ;</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Score+vie :</w:t></w:r></w:p><w:p><w:r><w:t>var vieX = [</w:t></w:r><w:r><w:t>108,124,140</w:t></w:r><w:r><w:t>] ;</w:t></w:r><w:r><w:br/><w:t>var vieY = [</w:t></w:r><w:r><w:t>341,341,341</w:t></w:r><w:r><w:t>] ;</w:t></w:r><w:r><w:br/><w:t>var vieL = [</w:t></w:r><w:r><w:t>15,15,15,15</w:t></w:r><w:r><w:t>] ;</w:t></w:r><w:r><w:br/><w:t>var vieH = [</w:t></w:r><w:r><w:t>19,19,19,19</w:t></w:r><w:r><w:t>] ;</w:t></w:r></w:p><w:p><w:r><w:t>ScoreX= [188,203,219,235</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>]</w:t></w:r><w:r><w:br/><w:t>ScoreY= 339</w:t></w:r><w:r><w:br/><w:t>ScoreL=11</w:t></w:r><w:r><w:br/><w:t>ScoreH=19</w:t></w:r></w:p><w:p/><w:p/>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$prange.InsertXML($xml)

Write-Output "Inserted score/vie blocks"
